$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(20081600, 0),
    @(20081700, 12801000000),
    @(20081800, 1732000000),
    @(20081900, 8825000000),
    @(20082000, 1225000000)
)

$startRow = 383
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
